# fix(publipostage): Correct status name
# Replace the "bleu" status label with "noir" and correct the status
# name wording from "pas de résultat ni de publication" to
# "pas de résultat postés ni publiés" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $labelCell = $ws.Cells.Item($r, 2)   # column B: statut_label
    $nameCell  = $ws.Cells.Item($r, 3)   # column C: statut_name

    if ($labelCell.Text -eq "bleu") {
        $labelCell.Value = "noir"
    }

    if ($nameCell.Text -eq "pas de résultat ni de publication") {
        $nameCell.Value = "pas de résultat postés ni publiés"
    }
}
